# Apply the edit: column N formulas reference column I instead of column K,
# and the sheet selection changes from N2:N1537 to Q5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2 is a standalone formula cell (not part of a shared formula group).
$ws.Range("N2").Formula = "=IF((MOD(ROW(I2)-2,3)=0), I2, 9999)"

# N3:N1536 covers all of the original shared-formula groups
# (N3:N66, N67:N130, ... , N1475:N1536). Assigning the whole block at once
# with a single relative formula preserves each of those original group
# boundaries/ids while swapping the K-column references for I-column ones.
$ws.Range("N3:N1536").Formula = "=IF((MOD(ROW(I3)-2,3)=0), I3, 9999)"

# N1537 was its own standalone formula cell before the edit too.
$ws.Range("N1537").Formula = "=IF((MOD(ROW(I1537)-2,3)=0), I1537, 9999)"

# Update the sheet's selection/active cell.
$ws.Range("Q5").Select()
